# Update decorator run time
$wb = $excel.ActiveWorkbook

# --- "Version" sheet: Start Time / End Time / Run Time strings ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = "09:35:01 03-22-2017 Wednesday Pacific Daylight Time"
$wsVersion.Range("C11").Value = "09:35:36 03-22-2017 Wednesday Pacific Daylight Time"
$wsVersion.Range("C12").Value = "0:0:34 0 days"

# --- "3_3" sheet: Deep Sleep (row3) / BT Idle (row4) stats + raw data ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.574
$ws33.Range("D3").Value = 0.604
$ws33.Range("E3").Value = 0.559
$ws33.Range("F3").Value = 0.011
$ws33.Range("H3").Value = "0.568721,0.572994,0.571051,0.573653,0.568816,0.570114,0.60046,0.565657,0.563993,0.579489,0.574672,0.563334,0.566683,0.571595,0.573639,0.563068,0.574006,0.571513,0.568075,0.566961,0.591194,0.573728,0.56919,0.589434,0.570304,0.569468,0.563809,0.57207,0.569828,0.579027,0.585351,0.569095,0.570671,0.57012,0.580609,0.566031,0.567892,0.573551,0.600005,0.568347,0.574944,0.567797,0.586934,0.568626,0.567879,0.565664,0.568911,0.590446,0.573273,0.573735,0.56743,0.565195,0.569557,0.572172,0.572253,0.558694,0.57745,0.592023,0.57671,0.564278,0.603823,0.566132,0.571506,0.573735,0.5715,0.565657,0.604468,0.561391,0.570766,0.569095,0.599815,0.572532,0.568816,0.597036,0.566316,0.566771,0.559346,0.594251,0.569095,0.573082,0.596764,0.574393,0.571975,0.564828,0.599828,0.575786,0.569747,0.572253,0.571323,0.570949,0.588598,0.572532,0.564461,0.603544,0.577552,0.568911,0.575317,0.57383,0.561853,0.562593"

$ws33.Range("C4").Value = 9.787000000000001
$ws33.Range("D4").Value = 9.851000000000001
$ws33.Range("E4").Value = 9.76
$ws33.Range("F4").Value = 0.021
$ws33.Range("H4").Value = "9.850566,9.843888,9.834602,9.830886,9.832006,9.846945,9.824853,9.810369,9.819099,9.838127,9.823739,9.812876,9.803684,9.814744,9.818453,9.813528,9.798019,9.800247,9.806381,9.810559,9.796728,9.780845,9.793657,9.804242,9.802285,9.777971,9.78168,9.787903,9.796538,9.782339,9.773426,9.783827,9.7906,9.787985,9.778528,9.762468,9.777224,9.784934,9.786606,9.775553,9.770362,9.782516,9.791714,9.773236,9.766089,9.768963,9.780947,9.788644,9.773046,9.766273,9.778712,9.784276,9.777129,9.775654,9.770913,9.787897,9.788922,9.779452,9.76219,9.77435,9.777883,9.783447,9.768318,9.766551,9.781035,9.788182,9.786701,9.773977,9.768318,9.777135,9.780947,9.776293,9.765431,9.764887,9.781775,9.788087,9.760899,9.773236,9.77649,9.791803,9.781694,9.765159,9.767659,9.784384,9.782238,9.776864,9.761354,9.77204,9.783535,9.787346,9.765254,9.76024,9.77825,9.785396,9.774812,9.760329,9.776673,9.79543,9.776952,9.776116"

# --- "1_8" sheet: Deep Sleep (row3) / BT Idle (row4) stats + raw data ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = -0
$ws18.Range("D3").Value = 0.008
$ws18.Range("E3").Value = -0.011
$ws18.Range("H3").Value = "0.008074,-0.001619,-0.001354,0.00533,0.003092,0.003936,-0.003574,0.000853,0.00311,0.001439,-0.000509,-0.002495,-0.000787,-0.000257,0.002561,-0.004702,-0.000226,0.002802,0.002265,0.000841,0.003079,-0.002476,0.002253,0.001433,-0.003556,-0.004129,-0.010789,-0.003568,-0.003013,2.7e-05,0.002512,-0.004955,-0.000528,-0.001367,0.002228,-0.001046,-0.001891,0.001162,-0.003839,-0.002464,0.00253,-0.000226,-0.001903,0.001963,0.004208,0.001414,-0.003556,0.001975,0.000304,0.002228,-0.002187,-0.000818,0.002247,6.4e-05,-0.001595,-0.005251,-0.002994,-0.000232,0.001143,0.002278,-0.002464,0.001969,-0.000497,-0.001058,-0.003297,-0.002199,0.000594,0.000304,-0.004696,-0.002181,0.002826,-0.003882,-0.003309,-0.007995,-0.004407,-0.001613,-0.001607,-0.00247,-0.003007,-0.004653,0.003634,-0.003303,0.001698,-0.003007,0.004744,0.001439,0.005034,-0.00247,-0.000528,0.006125,0.001433,0.004214,0.00504,0.001149,-0.003852,-0.002735,0.005873,-0.004943,0.00475,0.001402"

$ws18.Range("D4").Value = 0.008
$ws18.Range("E4").Value = -0.008
$ws18.Range("H4").Value = "-0.002705,-0.00662,-0.000793,0.003061,-0.002181,-0.001348,0.000323,-0.005541,0.004177,-0.001613,-0.003562,0.00171,0.00475,-0.001632,0.001149,-0.000269,-0.003007,-0.002754,-0.004974,-0.002748,0.00282,0.002789,-0.005498,0.004473,-0.002181,-0.002735,-0.001607,-0.003876,-0.001058,-0.003284,-0.001095,-0.003864,0.001717,-0.000257,0.000329,-0.004937,-0.003284,-0.002476,-0.002187,0.00253,-0.001342,-0.007718,-0.001064,-0.004974,0.002814,-0.000516,-0.003309,-0.003852,-0.006361,-0.000263,0.000329,0.002012,0.000853,0.002808,0.002537,-0.001891,-0.001619,-0.00136,-0.003025,0.000317,0.003369,0.00171,-0.00215,-0.001582,-0.001613,0.007556,-0.001058,0.00504,0.001445,0.001704,-0.001638,0.001717,-0.003013,-0.004684,-0.005492,0.002537,0.00171,0.001994,-0.000534,0.001143,0.003098,0.003098,-0.003833,-0.000232,0.00311,-0.003001,-0.006904,0.000582,-0.000799,-0.005794,-0.001354,-0.007422,-0.005769,-0.004413,-0.005547,3.3e-05,-0.001891,2.7e-05,-0.002211,-0.004955"
